# Generate Report for Handback
# Updates the localization-status workbook to reflect that the handback
# files are now in sync (both zh-cn and de-de), refreshes the "Latest
# Handback DateTime" timestamps, clears the stale "Error Detail" messages,
# and widens a couple of date/status columns to fit the new content.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": languages are now handed back & in sync with en-US
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------
# Sheet "zh-cn": refresh handback datetime, clear stale error detail
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("K2").Value = "2016-08-26 12:49:53"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.1
$zhcn.Columns.Item(16).ColumnWidth = 12.8

# ---------------------------------------------------------------------
# Sheet "de-de": refresh handback datetime, clear stale error detail
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K2").Value = "2016-08-26 12:50:03"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.1
$dede.Columns.Item(16).ColumnWidth = 12.8
